# Auto-generated edit script applying the policy_compteur.xlsx diff
# Updates cell values and reuses existing fill styles (s=2..5) without creating new styles

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells that already carry each target style (left untouched by this edit)
$styleSource = @{
    "1" = "B3"
    "2" = "B4"
    "3" = "L4"
    "4" = "D7"
    "5" = "H33"
}

function Set-CellStyleValue {
    param($Ref, $StyleIdx, $Value)
    $srcAddr = $styleSource[$StyleIdx]
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($Ref)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $dst.Value2 = $Value
}

$excel.CutCopyMode = $false

Set-CellStyleValue "E4" "3" 0
Set-CellStyleValue "K4" "2" 1
Set-CellStyleValue "D5" "2" 1
Set-CellStyleValue "E5" "2" 1
Set-CellStyleValue "G5" "2" 1
Set-CellStyleValue "H5" "2" 1
Set-CellStyleValue "I5" "2" 1
Set-CellStyleValue "D6" "2" 1
Set-CellStyleValue "E6" "2" 1
Set-CellStyleValue "H6" "2" 1
Set-CellStyleValue "G7" "2" 1
Set-CellStyleValue "H7" "2" 1
Set-CellStyleValue "D8" "2" 1
Set-CellStyleValue "E8" "2" 1
Set-CellStyleValue "F8" "2" 1
Set-CellStyleValue "G8" "2" 1
Set-CellStyleValue "H8" "2" 1
Set-CellStyleValue "I8" "2" 1
Set-CellStyleValue "D9" "2" 1
Set-CellStyleValue "F9" "2" 1
Set-CellStyleValue "G9" "2" 1
Set-CellStyleValue "H9" "2" 1
Set-CellStyleValue "I9" "2" 1
Set-CellStyleValue "E10" "2" 1
Set-CellStyleValue "E11" "2" 1
Set-CellStyleValue "K11" "3" 0
Set-CellStyleValue "C12" "3" 0
Set-CellStyleValue "D12" "2" 1
Set-CellStyleValue "B13" "2" 1
Set-CellStyleValue "D13" "3" 0
Set-CellStyleValue "E13" "2" 1
Set-CellStyleValue "G18" "2" 1
Set-CellStyleValue "B21" "2" 1
Set-CellStyleValue "C22" "2" 1
Set-CellStyleValue "E22" "2" 1
Set-CellStyleValue "G22" "3" 0
Set-CellStyleValue "G23" "3" 0
Set-CellStyleValue "G27" "3" 0
Set-CellStyleValue "B32" "5" 3
Set-CellStyleValue "F33" "2" 1
Set-CellStyleValue "G33" "5" 3
Set-CellStyleValue "F34" "2" 1
Set-CellStyleValue "B35" "2" 1
Set-CellStyleValue "C35" "3" 0
Set-CellStyleValue "F35" "2" 1
Set-CellStyleValue "G35" "5" 3
Set-CellStyleValue "H35" "5" 3
Set-CellStyleValue "B36" "5" 3
Set-CellStyleValue "D36" "5" 3
Set-CellStyleValue "G36" "5" 3
Set-CellStyleValue "B37" "2" 1
Set-CellStyleValue "C37" "2" 1
Set-CellStyleValue "D37" "3" 0
Set-CellStyleValue "F37" "2" 1
Set-CellStyleValue "G37" "2" 1
Set-CellStyleValue "J37" "3" 0
Set-CellStyleValue "F38" "2" 1
Set-CellStyleValue "D39" "2" 1
Set-CellStyleValue "E39" "3" 0
Set-CellStyleValue "F39" "2" 1
Set-CellStyleValue "I40" "5" 3
Set-CellStyleValue "J40" "3" 0
Set-CellStyleValue "F41" "3" 0
Set-CellStyleValue "I41" "5" 3

$excel.CutCopyMode = $false
